# "reserve, extend and history"
# Book 1 (row 2) has its loan extended by ~1h49m (both the lent date and
# the return date shift forward by the same amount), and is no longer
# shown as reserved (its "Reserved until" timestamp is cleared).
# Book 2 (row 3) also has its "Reserved until" history cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 45805.74464392842
$ws.Range("J2").Value = 45835.74464392842

$ws.Range("M2").Clear()
$ws.Range("M3").Clear()
